# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the row data for "Almeria" (row 47) and "Lugo" (row 48):
# the province names move, and the "Casos activos" (column C) values move
# together with them, so the correct figure stays paired with its province.
$ws.Range("A47").Value = "Lugo"
$ws.Range("C47").Value = 5

$ws.Range("A48").Value = "Almeria"
$ws.Range("C48").Value = 72

# Update the "last updated" timestamp string in cell A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 19:46"
